# Update the ANP Hamburger Marketshare model:
#  - Replace the instructional note text used throughout the "pairwise_comp"
#    sheet (previously: "Enter judgments for the paiwise comparisons in the
#    matrix or direct values in the green cells") with the new wording.
#  - Update the sheet's active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$note = "Enter pairwise comparisons in the white cells of the table or numerical data in the green cells. For the Direct Values column, if the smallest value is best, invert the value before entering it (e.g., `$10 as =1/10) ."

$noteCells = @(
    "A2", "A8", "A15", "A22", "A35", "A41", "A48", "A55", "A68", "A74",
    "A81", "A88", "A101", "A108", "A114", "A122", "A129", "A135", "A142",
    "A149", "A155", "A163", "A171", "A179", "A186", "A193", "A200", "A206",
    "A212", "A219", "A227", "A234", "A243", "A250", "A258", "A265", "A273",
    "A280", "A287", "A294", "A303", "A310", "A318", "A325", "A332", "A339"
)

foreach ($cellRef in $noteCells) {
    $ws.Range($cellRef).Value = $note
}

$null = $ws.Activate()
$null = $ws.Range("A339").Select()
$excel.ActiveWindow.ScrollRow = 320
$excel.ActiveWindow.ScrollColumn = 1
